$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-row updates: only D (Price) / E (Volume 1h) text changes ---
# Row -> @(NewD_or_$null, NewE_or_$null)
$updates = @{
    2  = @("30.739.78", "  +1.52%  ")
    3  = @("2.105.58",  "  +5.06%  ")
    4  = @("1.002",     "  +0.10%  ")
    5  = @("333.67",    "  +2.93%  ")
    6  = @($null,       "  +0.19%  ")
    7  = @("0.5290",    "  +3.66%  ")
    8  = @("0.4358",    "  +5.10%  ")
    9  = @("0.08953",   "  +2.83%  ")
    10 = @("47.10",     "  +10.21%  ")
    11 = @($null,       "  +2.53%  ")
    12 = @("24.73",     "  -0.94%  ")
    13 = @("2.106.38",  "  +5.32%  ")
    14 = @("6.711",     "  +2.68%  ")
    15 = @("7.753",     "  +4.25%  ")
    16 = @("96.53",     "  +2.47%  ")
    17 = @("1.002",     "  -0.08%  ")
    18 = @("0.00001131","  +1.20%  ")
    19 = @("0.06686",   "  +2.15%  ")
    20 = @("18.98",     "  -0.04%  ")
    21 = @("1.001",     "  +0.13%  ")
    22 = @($null,       "  +2.44%  ")
    23 = @("30.801.37", "  +1.53%  ")
    24 = @("12.26",     "  +4.64%  ")

    27 = @("22.53",     "  -0.28%  ")
    28 = @("2.565",     "  +6.32%  ")
    29 = @("162.55",    "  -0.73%  ")
    30 = @("132.67",    $null)
    31 = @("1.192",     "  +4.43%  ")
    32 = @("0.1079",    "  +2.58%  ")
    33 = @("6.162",     "  +1.51%  ")
    34 = @("3.946",     "  +2.96%  ")
    35 = @("1.539",     "  +14.77%  ")
    36 = @("0.02598",   "  +3.80%  ")
    37 = @("9.564",     "  +7.10%  ")
    38 = @("5.527",     "  +2.55%  ")
    39 = @("0.06741",   "  +2.30%  ")
    40 = @("12.63",     "  +2.68%  ")
    41 = @("0.2268",    "  +3.08%  ")
    42 = @("0.6810",    "  +2.67%  ")
    43 = @("1.243",     "  +1.25%  ")
    44 = @("1.001",     "  +0.19%  ")

    47 = @("2.217",     "  +0.54%  ")
    48 = @("3.657",     "  -0.13%  ")
    49 = @($null,       "  -0.35%  ")
    50 = @($null,       "  +3.57%  ")
    51 = @("1.193",     "  +8.07%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}

# --- Rows 25 & 26: coin identities swap (B/C), plus independent D/E updates ---
$ws.Cells.Item(25, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(25, 4).Value = "2.354.06"
$ws.Cells.Item(25, 5).Value = "  +5.38%  "

$ws.Cells.Item(26, 2).Value = "Toncoin"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(26, 4).Value = "2.285"
$ws.Cells.Item(26, 5).Value = "  +3.46%  "

# --- Rows 45 & 46: coin identities swap (B/C), plus independent D/E updates ---
$ws.Cells.Item(45, 2).Value = "Decentraland"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(45, 4).Value = "0.6405"
$ws.Cells.Item(45, 5).Value = "  +3.97%  "

$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).Value = "13.99"
$ws.Cells.Item(46, 5).Value = "  +2.75%  "
